$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the "Security Properties"/"Done" row down to row 3
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2

# Set the new row 2 values
$ws.Range("B2").Value = "Helloworld"
$ws.Range("C2").Value = "Done"

# Update the active selection as reflected in the diff
$ws.Range("G5").Select()
